# Modifications pour utiliser XGBClassifier et ajuster les prédictions
#
# Sheet "Valeurs réelles": rename the S+1/S+2/S+3 headers to *_class and
# replace the per-date price columns (C,D,E) with integer class labels.
# Sheet "Prédictions": replace the float PRED_S1/S2/S3 predictions with
# integer class predictions.

$wb = $excel.ActiveWorkbook

$wsReal = $wb.Worksheets.Item(1)   # "Valeurs réelles"
$wsPred = $wb.Worksheets.Item(2)   # "Prédictions"

# --- Sheet 1: "Valeurs réelles" ---------------------------------------

# Header renames: append "_class" to the S+1 / S+2 / S+3 headers.
$wsReal.Range("C1").Value = "PRIX EXP POMME FUJI FRANCE 170/220G CAT.I PLATEAU 1RG_S+1_class"
$wsReal.Range("D1").Value = "PRIX EXP POMME FUJI FRANCE 170/220G CAT.I PLATEAU 1RG_S+2_class"
$wsReal.Range("E1").Value = "PRIX EXP POMME FUJI FRANCE 170/220G CAT.I PLATEAU 1RG_S+3_class"

# Row-by-row replacement of columns C (S+1), D (S+2), E (S+3) with the new
# integer class values (row index => (C, D, E)).
$realClasses = @{
    2  = @(4, 2, 2)
    3  = @(2, 2, 2)
    4  = @(2, 2, 2)
    5  = @(2, 2, 2)
    6  = @(2, 2, 2)
    7  = @(2, 2, 2)
    8  = @(2, 2, 2)
    9  = @(2, 2, 2)
    10 = @(2, 2, 2)
    11 = @(2, 2, 2)
    12 = @(2, 2, 2)
    13 = @(2, 2, 2)
    14 = @(2, 2, 4)
    15 = @(2, 4, 2)
    16 = @(4, 2, 2)
    17 = @(2, 2, 4)
    18 = @(2, 4, 1)
    19 = @(4, 1, 2)
    20 = @(1, 2, 2)
    21 = @(2, 2, 2)
    22 = @(2, 2, 1)
    23 = @(2, 1, 2)
    24 = @(1, 2, 3)
    25 = @(2, 3, 1)
    26 = @(3, 1, 2)
    27 = @(1, 2, 2)
    28 = @(2, 2, 2)
}

foreach ($row in $realClasses.Keys) {
    $vals = $realClasses[$row]
    $wsReal.Cells.Item($row, 3).Value = $vals[0]   # C
    $wsReal.Cells.Item($row, 4).Value = $vals[1]   # D
    $wsReal.Cells.Item($row, 5).Value = $vals[2]   # E
}

# --- Sheet 2: "Prédictions" -------------------------------------------

# Row-by-row replacement of columns B (PRED_S1), C (PRED_S2), D (PRED_S3)
# with the new integer class predictions (row index => (B, C, D)).
$predClasses = @{
    2  = @(2, 2, 2)
    3  = @(0, 0, 0)
    4  = @(0, 0, 0)
    5  = @(0, 0, 0)
    6  = @(0, 0, 0)
    7  = @(0, 0, 0)
    8  = @(0, 0, 0)
    9  = @(0, 0, 0)
    10 = @(0, 0, 0)
    11 = @(0, 0, 0)
    12 = @(0, 0, 0)
    13 = @(0, -2, 0)
    14 = @(0, 0, 0)
    15 = @(0, 0, 0)
    16 = @(0, 0, 0)
    17 = @(-2, -2, 2)
    18 = @(2, 0, 0)
    19 = @(-2, 0, 0)
    20 = @(0, 0, 0)
    21 = @(0, 0, 0)
    22 = @(-2, -1, 0)
    23 = @(2, 2, 2)
    24 = @(2, 2, 2)
    25 = @(0, -2, 2)
    26 = @(0, 0, 0)
    27 = @(0, -2, 0)
    28 = @(0, 0, 0)
}

foreach ($row in $predClasses.Keys) {
    $vals = $predClasses[$row]
    $wsPred.Cells.Item($row, 2).Value = $vals[0]   # B
    $wsPred.Cells.Item($row, 3).Value = $vals[1]   # C
    $wsPred.Cells.Item($row, 4).Value = $vals[2]   # D
}
